$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated crypto price/volume figures (refreshed snapshot).
# Each entry: cell address, new text, and whether the text must be
# force-typed as Text (Excel would otherwise auto-convert a plain
# numeric-looking string like "1.006" into a Number and drop the
# trailing zero / "." thousands formatting).
$updates = @(
    @{ Cell = "D2"; Value = "28.379.70"; ForceText = $false }
    @{ Cell = "E2"; Value = "  -0.33%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "1.822.20"; ForceText = $false }
    @{ Cell = "E3"; Value = "  -0.37%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "1.006"; ForceText = $true }
    @{ Cell = "E4"; Value = "  +0.44%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "316.54"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +0.45%  "; ForceText = $false }
    @{ Cell = "E6"; Value = "  +0.41%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.5138"; ForceText = $true }
    @{ Cell = "E7"; Value = "  +0.30%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.3838"; ForceText = $true }
    @{ Cell = "E8"; Value = "  -2.00%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.08120"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +5.80%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "1.116"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +0.42%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "41.84"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +0.08%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "6.389"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +1.70%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "21.00"; ForceText = $true }
    @{ Cell = "E13"; Value = "  -0.36%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "1.005"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +0.41%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "7.424"; ForceText = $true }
    @{ Cell = "E15"; Value = "  -1.50%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "1.820.97"; ForceText = $false }
    @{ Cell = "E16"; Value = "  -0.33%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "93.85"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +0.14%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "0.00001104"; ForceText = $true }
    @{ Cell = "E18"; Value = "  +0.18%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "0.06636"; ForceText = $true }
    @{ Cell = "E19"; Value = "  -1.27%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "17.71"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +0.00%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "1.005"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +0.39%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "6.025"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -2.25%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "28.417.67"; ForceText = $false }
    @{ Cell = "E23"; Value = "  -0.29%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "11.45"; ForceText = $true }
    @{ Cell = "E24"; Value = "  +2.47%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "2.249"; ForceText = $true }
    @{ Cell = "E25"; Value = "  -0.31%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "159.78"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +2.03%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "20.96"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +1.39%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "2.027.39"; ForceText = $false }
    @{ Cell = "E28"; Value = "  -0.53%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "2.392"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -0.19%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "124.78"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +0.31%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "0.1095"; ForceText = $true }
    @{ Cell = "E31"; Value = "  +0.40%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "1.082"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -3.10%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "5.688"; ForceText = $true }
    @{ Cell = "E33"; Value = "  +0.33%  "; ForceText = $false }
    @{ Cell = "E34"; Value = "  +0.76%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "0.07399"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +5.26%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "12.28"; ForceText = $true }
    @{ Cell = "E36"; Value = "  +9.45%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "0.2200"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -0.67%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.02344"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +0.87%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "5.149"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -0.22%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "8.699"; ForceText = $true }
    @{ Cell = "E40"; Value = "  -2.99%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "0.6330"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +1.02%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "1.182"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +0.08%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "1.383"; ForceText = $true }
    @{ Cell = "E43"; Value = "  -0.67%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "13.53"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +0.71%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "0.6132"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +3.87%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "3.793"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +2.08%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "127.27"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +1.82%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "1.991"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +0.55%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "1.199"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +0.02%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "0.06902"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -0.33%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "1.068"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -0.08%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Temporarily switch to Text format so the numeric-looking string
        # is stored verbatim, then restore the original General format/style
        # so no visible formatting changes are left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

$wb.Save()